$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-8) are rotated/swapped. Columns D, M, N, O, P, S change values
# for rows 2,3,4,5 (cyclic shift: 2<-3, 3<-4, 4<-5, 5<-2) and rows 7,8 (swap: 7<-8, 8<-7).
# Row 6 is unchanged.

# Capture the "before" values for the affected rows/columns first, since we will
# overwrite them in place.
$cols = @("D", "M", "N", "O", "P", "S")
$rows = @(2, 3, 4, 5, 7, 8)

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row values, sourced from the original values of other rows (per the diff).
$mapping = @{
    2 = 3
    3 = 4
    4 = 5
    5 = 2
    7 = 8
    8 = 7
}

foreach ($r in $rows) {
    $srcRow = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig[$srcRow][$c]
    }
}
